$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the values that used to be on row 8 (the Oct 10 2020
# match vs Chennai Super Kings), replacing the Oct 25 2020 match.
$ws.Range("A2").Value = " Oct 10 2020"
$ws.Range("B2").Value = " Dubai (DSC)"
$ws.Range("C2").Value = "RCB won by 37 runs"
$ws.Range("D2").Value = "Royal Challengers Bangalore"
$ws.Range("E2").Value = "Chennai Super Kings"
$ws.Range("F2").Value = "Aaron Finch" + [char]160

# G2:K2 hold numeric-looking text ("2", "9", "0", "0", "22.22") that must
# stay text, not become real numbers - force text format before writing.
$ws.Range("G2:K2").NumberFormat = "@"
$ws.Range("G2").Value = "2"
$ws.Range("H2").Value = "9"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "22.22"

# Remove the now-duplicated rows 3 through 8 (the old Oct 12 / Sep 28 /
# Sep 21 / Oct 15 / Nov 6 / Oct 10 match rows).
$ws.Range("A3:K8").EntireRow.Delete()
